$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.228.53"
$ws.Range("E2").Value = "  +0.39%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.112.66"
$ws.Range("E3").Value = "  +0.36%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.69"
$ws.Range("E5").Value = "  +0.07%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.51"
$ws.Range("E6").Value = "  +0.64%  "

$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.520"
$ws.Range("E8").Value = "  -0.49%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.54"
$ws.Range("E9").Value = "  +1.55%  "

$ws.Range("E10").Value = "  -0.52%  "

$ws.Range("E11").Value = "  -0.16%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000248"
$ws.Range("E12").Value = "  -0.36%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "36.87"
$ws.Range("E13").Value = "  -1.13%  "

$ws.Range("E14").Value = "  -1.54%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.628.95"
$ws.Range("E15").Value = "  +0.41%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.211.32"
$ws.Range("E16").Value = "  +0.33%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.10"
$ws.Range("E17").Value = "  -1.31%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.114.59"
$ws.Range("E18").Value = "  +0.37%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.61"
$ws.Range("E19").Value = "  +1.95%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "493.00"
$ws.Range("E20").Value = "  +1.39%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.90"
$ws.Range("E21").Value = "  +4.67%  "

$ws.Range("B22").Value = "Polygon"
$ws.Range("C22").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.706"
$ws.Range("E22").Value = "  -1.67%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "84.01"
$ws.Range("E23").Value = "  -0.49%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.19"
$ws.Range("E24").Value = "  -1.68%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.30"
$ws.Range("E25").Value = "  -2.29%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.62"
$ws.Range("E26").Value = "  +6.21%  "

$ws.Range("E27").Value = "  -0.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.96"
$ws.Range("E28").Value = "  -1.29%  "

$ws.Range("E29").Value = "  -1.74%  "

$ws.Range("E30").Value = "  -0.23%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "28.43"
$ws.Range("E31").Value = "  -1.81%  "

$ws.Range("E32").Value = "  -0.18%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0₃0942"
$ws.Range("E33").Value = "  -6.05%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  -0.06%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.87"
$ws.Range("E35").Value = "  -0.64%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.974"
$ws.Range("E36").Value = "  -1.68%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "47.26"
$ws.Range("E37").Value = "  -0.70%  "

$ws.Range("E38").Value = "  -3.11%  "

$ws.Range("E39").Value = "  -1.64%  "

$ws.Range("E40").Value = "  +1.36%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.47"
$ws.Range("E41").Value = "  -2.25%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "389.12"
$ws.Range("E42").Value = "  +1.13%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.805.60"
$ws.Range("E43").Value = "  -1.12%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.58"
$ws.Range("E44").Value = "  -7.62%  "

$ws.Range("E45").Value = "  -2.52%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "134.99"
$ws.Range("E46").Value = "  -0.81%  "

$ws.Range("E47").Value = "  +0.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.08"
$ws.Range("E48").Value = "  +0.69%  "

$ws.Range("E49").Value = "  -0.72%  "

$ws.Range("E50").Value = "  -0.92%  "

$ws.Range("E51").Value = "  -1.89%  "
